$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 307
$firstCol = 19  # column S
$lastCol = 25   # column Y

for ($i = 2; $i -le $lastRow; $i++) {
    $a = $ws.Cells.Item($i, 1).Value2

    # Update "Förändrad" (Changed) date column C from 45184 to 45186
    $ws.Cells.Item($i, 3).Value = 45186

    for ($col = $firstCol; $col -le $lastCol; $col++) {
        $cell = $ws.Cells.Item($i, $col)
        if ($cell.HasFormula) {
            $f = $cell.Formula
            $newf = $f.Substring(0, $f.Length - 1) + ', "' + $a + '")'
            $cell.Formula = $newf
        }
    }
}
